$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.296.09"
$ws.Range("E2").Value = "  -15.34%  "
$ws.Range("D3").Value = "2.255.76"
$ws.Range("E3").Value = "  -22.27%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "431.36"
$ws.Range("E5").Value = "  -18.05%  "
$ws.Range("D6").Value = "116.54"
$ws.Range("E6").Value = "  -18.51%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.454"
$ws.Range("E8").Value = "  -16.91%  "
$ws.Range("D9").Value = "2.258.32"
$ws.Range("E9").Value = "  -22.43%  "
$ws.Range("D10").Value = "5.15"
$ws.Range("E10").Value = "  -13.78%  "
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  -21.28%  "
$ws.Range("D12").Value = "0.294"
$ws.Range("E12").Value = "  -18.05%  "
$ws.Range("D13").Value = "0.119"
$ws.Range("E13").Value = "  -7.23%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.646.23"
$ws.Range("E14").Value = "  -22.37%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "51.635.23"
$ws.Range("E15").Value = "  -14.79%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "18.29"
$ws.Range("E16").Value = "  -18.57%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000113"
$ws.Range("E17").Value = "  -18.92%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.273.43"
$ws.Range("E18").Value = "  -21.74%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "3.87"
$ws.Range("E19").Value = "  -22.12%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "289.91"
$ws.Range("E20").Value = "  -17.20%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "0.989"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "5.63"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "8.53"
$ws.Range("E23").Value = "  -26.39%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "4.98"
$ws.Range("E24").Value = "  -23.86%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "52.10"
$ws.Range("E26").Value = "  -19.45%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.385.53"
$ws.Range("E27").Value = "  -21.23%  "
$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").Value = "0.358"
$ws.Range("E28").Value = "  -20.66%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "0.138"
$ws.Range("E29").Value = "  -22.03%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "6.69"
$ws.Range("E31").Value = "  -14.21%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0625"
$ws.Range("E32").Value = "  -26.65%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "142.12"
$ws.Range("E33").Value = "  -6.61%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "16.38"
$ws.Range("E34").Value = "  -16.05%  "
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "1.29"
$ws.Range("E35").Value = "  -22.71%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "4.57"
$ws.Range("E36").Value = "  -17.59%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "0.780"
$ws.Range("E37").Value = "  -21.49%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "0.997"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "3.28"
$ws.Range("E39").Value = "  -23.51%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "0.965"
$ws.Range("E40").Value = "  -19.14%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "31.62"
$ws.Range("E41").Value = "  -15.93%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "10.12"
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.546"
$ws.Range("E43").Value = "  -15.92%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "0.0488"
$ws.Range("E44").Value = "  -15.64%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "3.03"
$ws.Range("E45").Value = "  -18.32%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.861.70"
$ws.Range("E46").Value = "  -18.64%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  -23.15%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.0801"
$ws.Range("E48").Value = "  -12.43%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0199"
$ws.Range("E49").Value = "  -16.14%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "3.89"
$ws.Range("E50").Value = "  -21.50%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "15.22"
$ws.Range("E51").Value = "  -25.24%  "
